$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 5808
$ws.Range("G3").Value = 8321
$ws.Range("K3").Value = 5977
$ws.Range("I4").Value = 1801
$ws.Range("J4").Value = 1833
$ws.Range("K4").Value = 1242
$ws.Range("K5").Value = 424
$ws.Range("K6").Value = 6570
$ws.Range("G7").Value = 24717
$ws.Range("I7").Value = 26259
$ws.Range("J7").Value = 29300
$ws.Range("K7").Value = 20021

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 581
$ws.Range("K8").Value = 1326
$ws.Range("K9").Value = 86
$ws.Range("K10").Value = 114
$ws.Range("K11").Value = 380
$ws.Range("K18").Value = 131
$ws.Range("K20").Value = 473
$ws.Range("K27").Value = 184
$ws.Range("K29").Value = 1093
$ws.Range("K33").Value = 869
$ws.Range("K37").Value = 675
$ws.Range("K42").Value = 739
$ws.Range("K48").Value = 255
$ws.Range("K49").Value = 109
$ws.Range("K50").Value = 97
$ws.Range("K51").Value = 256
$ws.Range("K52").Value = 526
$ws.Range("K53").Value = 253
$ws.Range("K54").Value = 386
$ws.Range("K55").Value = 222
$ws.Range("G63").Value = 291
$ws.Range("J63").Value = 115
$ws.Range("K64").Value = 129
$ws.Range("K65").Value = 465
$ws.Range("K67").Value = 780
$ws.Range("K72").Value = 94
$ws.Range("K76").Value = 271
$ws.Range("K77").Value = 139
$ws.Range("I79").Value = 748
$ws.Range("K79").Value = 497
$ws.Range("K84").Value = 157
$ws.Range("K85").Value = 939
$ws.Range("K88").Value = 214
$ws.Range("K89").Value = 293
$ws.Range("K90").Value = 184
$ws.Range("K91").Value = 227
$ws.Range("K92").Value = 77
$ws.Range("K95").Value = 335
$ws.Range("K99").Value = 332
$ws.Range("G101").Value = 24717
$ws.Range("I101").Value = 26259
$ws.Range("J101").Value = 29300
$ws.Range("K101").Value = 20021

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 188
$ws.Range("K7").Value = 581

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K6").Value = 123
$ws.Range("K7").Value = 380

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 90
$ws.Range("K7").Value = 293

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K6").Value = 234
$ws.Range("K7").Value = 939

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 142
$ws.Range("K3").Value = 155
$ws.Range("K7").Value = 526

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 65
$ws.Range("K3").Value = 66
$ws.Range("K7").Value = 253

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 368
$ws.Range("K3").Value = 404
$ws.Range("K6").Value = 443
$ws.Range("K7").Value = 1326

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 232
$ws.Range("K3").Value = 319
$ws.Range("K6").Value = 256
$ws.Range("K7").Value = 869

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 117
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K6").Value = 195
$ws.Range("K7").Value = 675

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 152
$ws.Range("K6").Value = 172
$ws.Range("K7").Value = 465

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 86
$ws.Range("K7").Value = 332

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 280
$ws.Range("K4").Value = 44
$ws.Range("K5").Value = 19
$ws.Range("K7").Value = 780

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K2").Value = 53
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K3").Value = 21
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 63
$ws.Range("K6").Value = 206
$ws.Range("K7").Value = 386

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 313
$ws.Range("K3").Value = 392
$ws.Range("K4").Value = 52
$ws.Range("K7").Value = 1093

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value = 43
$ws.Range("K3").Value = 45

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value = 51
$ws.Range("K6").Value = 140
$ws.Range("K7").Value = 271

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 200
$ws.Range("K6").Value = 274
$ws.Range("K7").Value = 739

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 222

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 109
$ws.Range("K7").Value = 227

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I4").Value = 39
$ws.Range("K6").Value = 123
$ws.Range("I7").Value = 748
$ws.Range("K7").Value = 497

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 154
$ws.Range("K4").Value = 20
$ws.Range("K5").Value = 8
$ws.Range("K6").Value = 135
$ws.Range("K7").Value = 473

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K4").Value = 11
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 214

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 50
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 67
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 73
$ws.Range("K4").Value = 29
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 256

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K3").Value = 56
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 139
